# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Reorders the two workers' mora-period detail rows (16-36): WILFREDO PEREZ
# HERAZO now appears first (periods 1811 -> 1801 descending) followed by
# JOSE DE JESUS NAVARRO (periods 1712 -> 1703 descending). The F (Valor
# Mora) / G (Salario Basico) pair for each worker+period combination is
# carried along with its row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# row => (DocNo, Name, Period, ValorMora, SalarioBasico)
$rows = @(
    @{ Row=16; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1811"; F=31249; G=737717 },
    @{ Row=17; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1810"; F=31249; G=737717 },
    @{ Row=18; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1809"; F=31249; G=737717 },
    @{ Row=19; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1808"; F=29509; G=737717 },
    @{ Row=20; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1807"; F=29509; G=737717 },
    @{ Row=21; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1806"; F=29509; G=737717 },
    @{ Row=22; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1805"; F=29509; G=737717 },
    @{ Row=23; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1804"; F=29509; G=737717 },
    @{ Row=24; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1803"; F=29509; G=737717 },
    @{ Row=25; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1802"; F=29509; G=737717 },
    @{ Row=26; Doc="9162606";    Name="WILFREDO PEREZ HERAZO"; Period="1801"; F=29509; G=737717 },
    @{ Row=27; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1712"; F=29509; G=781242 },
    @{ Row=28; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1711"; F=29509; G=781242 },
    @{ Row=29; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1710"; F=29509; G=781242 },
    @{ Row=30; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1709"; F=29509; G=781242 },
    @{ Row=31; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1708"; F=29509; G=781242 },
    @{ Row=32; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1707"; F=29509; G=781242 },
    @{ Row=33; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1706"; F=29509; G=781242 },
    @{ Row=34; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1705"; F=29509; G=781242 },
    @{ Row=35; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1704"; F=29509; G=781242 },
    @{ Row=36; Doc="9692332";    Name="JOSE DE JESUS NAVARRO"; Period="1703"; F=29509; G=781242 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc      # C: N° Doc Trabajador
    $ws.Cells.Item($r.Row, 4).Value = $r.Name     # D: Nombre Trabajador
    $ws.Cells.Item($r.Row, 5).Value = $r.Period   # E: Periodo Mora
    $ws.Cells.Item($r.Row, 6).Value = $r.F        # F: Valor Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.G        # G: Salario Basico
}
